$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Enter the new "Waypart" input values for column J (rows 3-6).
$ws.Range("J3").Value = 1596370
$ws.Range("J4").Value = 25285203
$ws.Range("J5").Value = 23560038
$ws.Range("J6").Value = 4760166

# Update the current selection on the sheet to match the recorded cursor position.
$ws.Range("M10").Select()
